# Append a new, blank slide to the end of the deck.
#
# The author's commit renamed the source file and did further work on the
# MEng flow-chart report; the only structural change that lands inside this
# particular presentation is a brand new (still empty) slide appended after
# the existing four slides -- the author's sketchpad for the new flow chart
# they mention getting "inspiration" for. It uses the same "Blank" layout
# as the other content slides in the deck (slides 2-4).

$p = $ppt.ActivePresentation

# ppLayoutBlank = 12; insert after the last existing slide (index 5 = end).
$newSlideIndex = $p.Slides.Count + 1
$newSlide = $p.Slides.Add($newSlideIndex, 12)
